$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Empty in nominal" mistake: the relation-type marker in row 2 was
# attached to the wrong column block. B2 (first block, columns B:J)
# should be 0 and K2 (second block, columns K:S) should be 1.
$ws.Range("B2").Value = 0
$ws.Range("K2").Value = 1

# Row 4 holds the summary statistics (AUC, KS, KS_p_val, rel_type, GINI,
# Count, Empty, Empty% in level, Empty% in all Empty) for each of the
# three column blocks (B:J, K:S, T:AB). Correcting the bug above changes
# the computed statistics for all three blocks.

# Block 1: B4:J4
$ws.Range("B4").Value = 0.7052419430788646
$ws.Range("C4").Value = 0.2992315081652258
$ws.Range("D4").Value = [double]"2.49770031077616E-17"
$ws.Range("E4").Value = -1
$ws.Range("F4").Value = 0.4104838861577291
$ws.Range("G4").Value = 306
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0

# Block 2: K4:S4
$ws.Range("K4").Value = 0.5644922286617492
$ws.Range("L4").Value = 0.1144296628029505
$ws.Range("M4").Value = 0.003151239620582965
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0.1289844573234984
$ws.Range("P4").Value = 584
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0

# Block 3: T4:AB4 (only the values that actually changed)
$ws.Range("T4").Value = 0.7851685393258426
$ws.Range("U4").Value = 0.4502553626149131
$ws.Range("V4").Value = [double]"1.186355338027967E-18"
$ws.Range("W4").Value = 1
$ws.Range("X4").Value = 0.5703370786516853
$ws.Range("Y4").Value = 110
$ws.Range("Z4").Value = 0
